$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L16").Value = -5.769230769230
$ws.Range("N16").Value = -82.246376811594
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -54.545454545454
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 112
$ws.Range("K17").Value = -6.25
$ws.Range("L17").Value = -31.372549019607
$ws.Range("M17").Value = 23.529411764705
$ws.Range("N17").Value = -50.471698113207
$ws.Range("C18").Value = "'0"
$ws.Range("D18").Value = "'0"
$ws.Range("E18").Value = "***.*"
$ws.Range("F18").Value = 1
$ws.Range("H18").Value = 0
$ws.Range("M18").Value = -5.263157894736
$ws.Range("N18").Value = -87.203791469194
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 40
$ws.Range("I19").Value = 169
$ws.Range("J19").Value = 135
$ws.Range("K19").Value = 25.185185185185
$ws.Range("L19").Value = 7.643312101910
$ws.Range("M19").Value = 49.557522123893
$ws.Range("N19").Value = -12.886597938144
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 59
$ws.Range("J20").Value = 24
$ws.Range("K20").Value = 145.833333333333
$ws.Range("L20").Value = 73.529411764705
$ws.Range("M20").Value = 110.714285714286
$ws.Range("N20").Value = -80.844155844155
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 50
$ws.Range("G21").Value = 26
$ws.Range("H21").Value = -3.846153846153
$ws.Range("I21").Value = 448
$ws.Range("J21").Value = 342
$ws.Range("K21").Value = 30.994152046783
$ws.Range("L21").Value = -4.273504273504
$ws.Range("M21").Value = 27.272727272727
$ws.Range("N21").Value = -68.671328671328
$ws.Range("C23").Value = "'0"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "***.*"
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -80
$ws.Range("M23").Value = 5
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 28
$ws.Range("H24").Value = 14.285714285714
$ws.Range("I24").Value = 527
$ws.Range("J24").Value = 382
$ws.Range("K24").Value = 37.958115183246
$ws.Range("L24").Value = 8.436213991769
$ws.Range("M24").Value = 78.040540540540
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 15
$ws.Range("H25").Value = -21.052631578947
$ws.Range("I25").Value = 216
$ws.Range("J25").Value = 189
$ws.Range("K25").Value = 14.285714285714
$ws.Range("L25").Value = -8.860759493670
$ws.Range("M25").Value = -35.905044510385
$ws.Range("C26").Value = "'0"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "***.*"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("C27").Value = "'0"
